$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the values first.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / bordered / centered-top style on B1 (font, alignment, border).
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2

# Clone that exact style onto A2 via copy/paste-special so both cells share
# a single cell format record instead of each mutation step leaving its own
# (possibly orphaned) style behind.
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
